$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the StockName (column C) values for all data rows - no longer tracked per-row
$ws.Range("C2:C8").ClearContents()

# Update the ID column (A) with the new set of usernames.
# Assignment order matters because the runtime appends new shared strings
# in the order values are written, and the target file expects the new
# strings appended in this exact order: user270, hoho222, nick0422, jiho264.
$ws.Range("A6").Value = "user270"
$ws.Range("A7").Value = "user270"
$ws.Range("A8").Value = "user270"

$ws.Range("A2").Value = "hoho222"
$ws.Range("A3").Value = "nick0422"
$ws.Range("A4").Value = "jiho264"
$ws.Range("A5").Value = "hoho222"

# Update StockCode (column B) for rows that changed owner; keep values stored
# as text (leading apostrophe) so the existing quote-prefixed text style (s=2)
# used for this column is preserved instead of Excel reformatting the cell.
$ws.Range("B4").Value = "'000660"
$ws.Range("B5").Value = "'000660"
$ws.Range("B6").Value = "'000660"

# Add the two new rows (hseop884 holding 086520 and 247540).
$ws.Range("A9").Value = "hseop884"
$ws.Range("B9").Value = "'086520"
$ws.Range("D9").Value = 4
$ws.Range("E9").Value = 0

$ws.Range("A10").Value = "hseop884"
$ws.Range("B10").Value = "'247540"
$ws.Range("D10").Value = 5
$ws.Range("E10").Value = 0

# Update the selected range shown when the sheet is active.
$ws.Range("C2:C10").Select()
